$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New row 23 - "79. Word Search" (Backtracking / Medium), added under the
# existing Backtracking block (rows 18-22).

$ws.Range("A23").Value = "Backtracking"
$ws.Range("A23").Style = $ws.Range("A22").Style

$ws.Range("B23").Value = "Medium"
$ws.Range("B23").Style = $ws.Range("B22").Style

$ws.Range("C23").Value = "79. Word Search"

$noteText = "We basically need to run a dfs with backtracking. First find the starting word in the matrix and then run dfs(rowIdx, colIdx, 0) search from that element and return True if dfs returns true as its final return value, otherwise return False at the end.`nBase case will be return True if check and match upto the last letter of word, or  False if ``the row or col idx is out of bounds, or the word does not match, or if its already visited.`nAdd the current element to a set or mark it in the matrix with ""#"" to signify we have visited it, then check the top, bottom, left and right element for the next word match recursively, then de mark it and return the result of the checks we ran.`n"
$ws.Range("D23").Value = $noteText

# Hyperlink for the problem name (added before re-applying the cell style,
# since Hyperlinks.Add stomps on the cell's existing style).
$ws.Hyperlinks.Add($ws.Range("C23"), "https://leetcode.com/problems/word-search/", [System.Type]::Missing, [System.Type]::Missing, "https://leetcode.com/problems/word-search/")

$ws.Range("C23").Style = $ws.Range("C22").Style
$ws.Range("D23").Style = $ws.Range("D22").Style

# Append the "Optimizations" rich-text run to the notes cell so the shared
# string carries the same multi-run formatting as the source workbook
# (bold "Optimizations" label followed by the regular continuation text).
$optRuns = $ws.Range("D23").GetRichTextRuns()
